# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    3  = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; G = 3.811642989160245 }
    4  = @{ B = 0.3464964993005633; C = 0.3375848360084654; D = 0.1529057820181812; E = 0.4998867070740569; G = 1.336873824401267 }
    5  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    6  = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    7  = @{ B = 0.7287194209349384; C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 5.964442013611383 }
    8  = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    9  = @{ B = 0.3464964993005633; C = 0.3375848360084654; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.896700893398075 }
    10 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    11 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; G = 3.811642989160245 }
    12 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    13 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
    14 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
